$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy E2's cell formatting (italic font style) down to E3:E5 so the whole
# out:cost column in the table shares the same look.
$ws.Range("E2").Copy()
$ws.Range("E3:E5").PasteSpecial(-4122)

# Fill in the remaining lab test-case data (miles, MPG, price) for the
# three rows that were still blank placeholders.
$ws.Range("B3").Value = 40
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 2

$ws.Range("B4").Value = 50
$ws.Range("C4").Value = 25
$ws.Range("D4").Value = 3

$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 2

# Replace the "insert Excel formula here" placeholder in E2, and extend the
# same cost formula down the whole out:cost column (E2:E5) in one shot so
# Excel records it as a shared formula.
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Item(5)
$col.DataBodyRange.Formula = "=(B2/C2*D2)"

# Match the selection left by the author's last save.
$ws.Range("E2").Select() | Out-Null
